$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the test case identifiers in column A:
# "test_InputMethod_SCB_func_..." -> "test_InputMethod_SCB_Func_..."
for ($r = 1; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val -like "*_func_*") {
        $cell.Value = $val -replace "_func_", "_Func_"
    }
}

# Move the active selection from H16 to B6
$ws.Range("B6").Select()
